# PowerSlate Integration Fields.docx - apply commit diff
# "Bugfixes, documentation, new tool"
#  - Fix Comments and County bugs
#  - Update fields doc
#  - New tool to select Religion mapping XML

$d = $word.ActiveDocument

$pkgNsAttr = 'xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"'
$wNsAttr   = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphWithXml($paragraph, [string]$bodyInnerXml) {
    # Replaces an entire paragraph (incl. its end-of-paragraph mark) with
    # the supplied <w:p>...</w:p> markup by collapsing to the end of the
    # paragraph's range and round-tripping through Range.InsertXML, which
    # Word resolves against the paragraph that the (zero-length) range sits
    # inside of.
    $insertRange = $paragraph.Range
    $insertRange.Collapse(0)
    $xml = '<?xml version="1.0" encoding="utf-8"?><pkg:package ' + $pkgNsAttr + '>' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
        '<w:document ' + $wNsAttr + '><w:body>' + $bodyInnerXml + '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $insertRange.InsertXML($xml)
}

function Append-ParagraphsAfter($paragraph, [string]$bodyInnerXml) {
    # Appends new sibling paragraph(s) right after $paragraph, inside the
    # same parent (table cell or body), WITHOUT disturbing $paragraph
    # itself. Works by placing a collapsed range right before the
    # paragraph mark (so it lands "inside" the existing paragraph's
    # container) and inserting full <w:p> markup there.
    $full = $paragraph.Range.Text
    $trimmed = $full.TrimEnd([char]13, [char]7)
    $insertPos = $paragraph.Range.Start + $trimmed.Length
    $insertRange = $d.Range($insertPos, $insertPos)
    $xml = '<?xml version="1.0" encoding="utf-8"?><pkg:package ' + $pkgNsAttr + '>' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
        '<w:document ' + $wNsAttr + '><w:body>' + $bodyInnerXml + '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $insertRange.InsertXML($xml)
}

function Find-ParagraphByText([string]$text) {
    $range = $d.Content
    $found = $range.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return $null }
    return $range.Paragraphs(1)
}

# ---------------------------------------------------------------
# 1. Document the new Stops integration fields (StopCode, StopDate,
#    Cleared, ClearedDate, Comments) under the "Stops" row.
# ---------------------------------------------------------------
$cnf = 'w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"'

$stopsXml = ('<w:p><w:pPr><w:cnfStyle {0}/></w:pPr></w:p>' -f $cnf) + `
    ('<w:p><w:pPr><w:cnfStyle {0}/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>StopCode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Code Value Key</w:t></w:r></w:p>' -f $cnf) + `
    ('<w:p><w:pPr><w:cnfStyle {0}/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>StopDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Date</w:t></w:r></w:p>' -f $cnf) + `
    ('<w:p><w:pPr><w:cnfStyle {0}/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Cleared</w:t></w:r><w:r><w:t>: Y/N</w:t></w:r></w:p>' -f $cnf) + `
    ('<w:p><w:pPr><w:cnfStyle {0}/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ClearedDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Date</w:t></w:r></w:p>' -f $cnf) + `
    ('<w:p><w:pPr><w:cnfStyle {0}/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Comments</w:t></w:r><w:r><w:t>: Free Text</w:t></w:r></w:p>' -f $cnf)

$stopsPara = Find-ParagraphByText "matched by stop code and stop date."
if ($stopsPara -ne $null) {
    Append-ParagraphsAfter $stopsPara $stopsXml
}

# ---------------------------------------------------------------
# 2. Mark a page break before "Students are looked up in PowerFAIDS..."
# ---------------------------------------------------------------
$studentsPara = Find-ParagraphByText "Students are looked up in PowerFAIDS via PCID and SSN (government ID)."
if ($studentsPara -ne $null) {
    $xml = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Students are looked up in PowerFAIDS via PCID and SSN (government ID).</w:t></w:r></w:p>'
    Replace-ParagraphWithXml $studentsPara $xml
}

# ---------------------------------------------------------------
# 3. Remove the (now stale) page break from the "Field/Data" header
#    cell -- it moved earlier in the document (see #2 above).
# ---------------------------------------------------------------
$fieldPara = Find-ParagraphByText "Field/Data"
if ($fieldPara -ne $null) {
    $xml = '<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Field</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>/Data</w:t></w:r></w:p>'
    Replace-ParagraphWithXml $fieldPara $xml
}

# ---------------------------------------------------------------
# 4. Bump the "Updated" date in the header from 2024-01-10 to 2024-04-08.
# ---------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdrPara = $hdr.Range.Paragraphs(1)
if ($hdrPara -ne $null) {
    $insertRange = $hdrPara.Range
    $insertRange.Collapse(0)
    $xml = '<?xml version="1.0" encoding="utf-8"?><pkg:package ' + $pkgNsAttr + '>' + `
        '<pkg:part pkg:name="/word/header1.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.header+xml"><pkg:xmlData>' + `
        '<w:hdr ' + $wNsAttr + '><w:p><w:pPr><w:pStyle w:val="Header"/><w:jc w:val="right"/></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">Updated </w:t></w:r>' + `
        '<w:r><w:t>2024-0</w:t></w:r>' + `
        '<w:r><w:t>4-08</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:t>by Wyatt Best</w:t></w:r>' + `
        '</w:p></w:hdr></pkg:xmlData></pkg:part></pkg:package>'
    $insertRange.InsertXML($xml)
}
